$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 398.5
$ws.Range("I9").Value = 432.33334
$ws.Range("J9").Value = 297
$ws.Range("K9").Value = 432.33334
$ws.Range("L9").Value = 297
$ws.Range("M9").Value = -263.33334
$ws.Range("N9").Value = -635

$ws.Range("H51").Value = 1000
$ws.Range("I51").Value = 1000
$ws.Range("K51").Value = 1000
$ws.Range("M51").Value = -516

$ws.Range("H62").Value = 50
$ws.Range("I62").Value = 50
$ws.Range("K62").Value = 50
$ws.Range("M62").Value = 574

$ws.Range("H65").Value = 50
$ws.Range("I65").Value = 50
$ws.Range("K65").Value = 250
$ws.Range("M65").Value = 2870

$ws.Range("H103").Value = 1399.8
$ws.Range("I103").Value = 1500
$ws.Range("K103").Value = 4500
$ws.Range("M103").Value = -3914

$ws.Range("H125").Value = 1645.9
$ws.Range("I125").Value = 820.6667
$ws.Range("J125").Value = 1999.5714
$ws.Range("K125").Value = 7386.0003
$ws.Range("L125").Value = 17996.1426
$ws.Range("M125").Value = -4926.0003
$ws.Range("N125").Value = -22916.1426

$ws.Range("H132").Value = 1843.8846
$ws.Range("I132").Value = 1877.64
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 5632.92
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -3102.92
$ws.Range("N132").Value = -8060

$ws.Range("H135").Value = 1875.2307
$ws.Range("I135").Value = 1875.2307
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 16877.0763
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 2051.1384
$ws.Range("I137").Value = 1519.3829
$ws.Range("J137").Value = 3439.611
$ws.Range("K137").Value = 4558.1487
$ws.Range("L137").Value = 10318.833
$ws.Range("M137").Value = -2008.1487
$ws.Range("N137").Value = -15418.833

$ws.Range("H138").Value = 5045.8667
$ws.Range("J138").Value = 5182.5835
$ws.Range("L138").Value = 15547.7505
$ws.Range("N138").Value = -25827.7505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2995.25
$ws.Range("I2").Value = 2994.5
$ws.Range("J2").Value = 2996
$ws.Range("K2").Value = 2994.5
$ws.Range("L2").Value = 2996
$ws.Range("M2").Value = -2881.5
$ws.Range("N2").Value = -3222

$ws.Range("H32").Value = 21332.084
$ws.Range("J32").Value = 24999
$ws.Range("L32").Value = 24999
$ws.Range("N32").Value = -25573

$ws.Range("H45").Value = 3998.5
$ws.Range("I45").Value = 3997.5
$ws.Range("J45").Value = 3999.5
$ws.Range("K45").Value = 3997.5
$ws.Range("L45").Value = 3999.5
$ws.Range("M45").Value = -3620.5
$ws.Range("N45").Value = -4753.5

$ws.Range("H61").Value = 1764.1034
$ws.Range("I61").Value = 1626.9131
$ws.Range("K61").Value = 1626.9131
$ws.Range("M61").Value = -1414.9131

$ws.Range("H74").Value = 2063.238
$ws.Range("I74").Value = 1341.2858
$ws.Range("K74").Value = 1341.2858
$ws.Range("M74").Value = -467.2858000000001

$ws.Range("H77").Value = 2063.238
$ws.Range("I77").Value = 1341.2858
$ws.Range("K77").Value = 6706.429
$ws.Range("M77").Value = -2338.429

$ws.Range("H97").Value = 223.375
$ws.Range("I97").Value = 169.57143
$ws.Range("K97").Value = 169.57143
$ws.Range("M97").Value = 326.42857

$ws.Range("H116").Value = 2995.25
$ws.Range("I116").Value = 2994.5
$ws.Range("J116").Value = 2996
$ws.Range("K116").Value = 2994.5
$ws.Range("L116").Value = 2996
$ws.Range("M116").Value = -700.5
$ws.Range("N116").Value = -7584

$ws.Range("H132").Value = 1392.1333
$ws.Range("I132").Value = 1392.4615
$ws.Range("K132").Value = 4177.3845
$ws.Range("M132").Value = -1647.3845

$ws.Range("H136").Value = 1764.1034
$ws.Range("I136").Value = 1626.9131
$ws.Range("K136").Value = 4880.7393
$ws.Range("M136").Value = -2330.7393

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2995.25
$ws.Range("I3").Value = 2994.5
$ws.Range("J3").Value = 2996
$ws.Range("K3").Value = 2994.5
$ws.Range("L3").Value = 2996
$ws.Range("M3").Value = -2880.5
$ws.Range("N3").Value = -3224

$ws.Range("H20").Value = 8201.5
$ws.Range("I20").Value = 6602.3335
$ws.Range("K20").Value = 6602.3335
$ws.Range("M20").Value = -6355.3335

$ws.Range("H134").Value = 4277
$ws.Range("I134").Value = 4485.5
$ws.Range("K134").Value = 13456.5
$ws.Range("M134").Value = -10921.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2489.1
$ws.Range("I31").Value = 2489.1
$ws.Range("K31").Value = 2489.1
$ws.Range("M31").Value = -2194.1

$ws.Range("H34").Value = 2489.1
$ws.Range("I34").Value = 2489.1
$ws.Range("K34").Value = 2489.1
$ws.Range("M34").Value = -2287.1

$ws.Range("H58").Value = 2884.125
$ws.Range("J58").Value = 2199.5
$ws.Range("L58").Value = 2199.5
$ws.Range("N58").Value = -2605.5

$ws.Range("H107").Value = 1441.1111
$ws.Range("J107").Value = 2228.1428
$ws.Range("L107").Value = 2228.1428
$ws.Range("N107").Value = -6068.1428

$ws.Range("H132").Value = 4778
$ws.Range("I132").Value = 4741.1665
$ws.Range("K132").Value = 14223.4995
$ws.Range("M132").Value = -11693.4995

$ws.Range("H136").Value = 2884.125
$ws.Range("J136").Value = 2199.5
$ws.Range("L136").Value = 6598.5
$ws.Range("N136").Value = -11698.5

$ws.Range("H141").Value = 57883.9
$ws.Range("J141").Value = 57883.9
$ws.Range("L141").Value = 57883.9
$ws.Range("N141").Value = -68243.89999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4080.1667
$ws.Range("I113").Value = 3995
$ws.Range("J113").Value = 4097.2
$ws.Range("K113").Value = 11985
$ws.Range("L113").Value = 12291.6
$ws.Range("M113").Value = -9815
$ws.Range("N113").Value = -16631.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2302.3333
$ws.Range("I102").Value = 2101
$ws.Range("K102").Value = 2101
$ws.Range("M102").Value = -479

$ws.Range("H132").Value = 1806.8
$ws.Range("I132").Value = 1806.8
$ws.Range("K132").Value = 5420.4
$ws.Range("M132").Value = -2890.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3159.7334
$ws.Range("I61").Value = 3145.923
$ws.Range("K61").Value = 3145.923
$ws.Range("M61").Value = -2943.923

$ws.Range("H93").Value = 2247.6843
$ws.Range("I93").Value = 2302.5625
$ws.Range("J93").Value = 1955
$ws.Range("K93").Value = 2302.5625
$ws.Range("L93").Value = 1955
$ws.Range("M93").Value = -1054.5625
$ws.Range("N93").Value = -4451

$ws.Range("H113").Value = 3159.7334
$ws.Range("I113").Value = 3145.923
$ws.Range("K113").Value = 3145.923
$ws.Range("M113").Value = -975.9229999999998

$ws.Range("H132").Value = 3814.1428
$ws.Range("I132").Value = 3400
$ws.Range("K132").Value = 10200
$ws.Range("M132").Value = -7670

$ws.Range("H136").Value = 7731.0835
$ws.Range("I136").Value = 7698.1
$ws.Range("K136").Value = 23094.3
$ws.Range("M136").Value = -20544.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 45000
$ws.Range("I61").Value = 45000
$ws.Range("K61").Value = 45000
$ws.Range("M61").Value = -44708

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H113").Value = 1836.1111
$ws.Range("J113").Value = 379.6
$ws.Range("L113").Value = 1138.8
$ws.Range("N113").Value = -5478.8

$ws.Range("H122").Value = 3233.9092
$ws.Range("I122").Value = 3186.6667
$ws.Range("K122").Value = 9560.000100000001
$ws.Range("M122").Value = -7110.000100000001

$ws.Range("H132").Value = 3360.875
$ws.Range("I132").Value = 2617.0908
$ws.Range("K132").Value = 7851.2724
$ws.Range("M132").Value = -5321.2724

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
